# Applies the "Add files via upload" refresh of Saldo_guide.xlsx:
#  - Sheet (tab) name date/time stamp bumped from 20241205-102334 to 20241206-093348
#  - Every row's "Dt. Referencia" (column G) bumped one day: 45631 -> 45632
#  - A handful of rows got updated "Vl. Projetado" (D) / "Saldo Previsto" (E)
#    figures, with "Vl. Total" (H) following as D+E for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to match the new export timestamp.
$ws.Name = "IClientBalance-20241206-093348-"

# Column G ("Dt. Referencia") moves from 45631 (2024-12-05) to 45632 (2024-12-06)
# for every data row (2 through 274).
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45632
}

# Row-specific updates to Vl. Projetado (D), Saldo Previsto (E) and the
# resulting Vl. Total (H = D + E).
$ws.Cells.Item(5, 4).Value = 41510.01
$ws.Cells.Item(5, 8).Value = 42391.76

$ws.Cells.Item(6, 5).Value = 188.68
$ws.Cells.Item(6, 8).Value = 188.68

$ws.Cells.Item(15, 4).Value = 164805.70000000001
$ws.Cells.Item(15, 8).Value = 165761.85999999999

$ws.Cells.Item(49, 4).Value = 5704.27
$ws.Cells.Item(49, 8).Value = 6394.39

$ws.Cells.Item(104, 4).Value = 112275.53
$ws.Cells.Item(104, 8).Value = 112698.47

$ws.Cells.Item(112, 5).Value = 974.5
$ws.Cells.Item(112, 8).Value = 974.5

$ws.Cells.Item(118, 5).Value = 707.16
$ws.Cells.Item(118, 8).Value = 707.16

$ws.Cells.Item(132, 4).Value = 35747.03
$ws.Cells.Item(132, 8).Value = 36535.89

$ws.Cells.Item(143, 4).Value = 164493.85999999999
$ws.Cells.Item(143, 8).Value = 207485.4

$ws.Cells.Item(231, 5).Value = 14408.65
$ws.Cells.Item(231, 8).Value = 14408.65

$ws.Cells.Item(235, 4).Value = 23827.74
$ws.Cells.Item(235, 8).Value = 24712.61

$ws.Cells.Item(264, 4).Value = 79038.149999999994
$ws.Cells.Item(264, 8).Value = 79379.55

$ws.Cells.Item(265, 4).Value = 42499.25
$ws.Cells.Item(265, 8).Value = 43473.7

$ws.Cells.Item(270, 4).Value = 30302.18
$ws.Cells.Item(270, 5).Value = 906.76
$ws.Cells.Item(270, 8).Value = 31208.94

$ws.Cells.Item(271, 4).Value = 41249.760000000002
$ws.Cells.Item(271, 8).Value = 42224.25

$ws.Cells.Item(273, 4).Value = 27469.87
$ws.Cells.Item(273, 8).Value = 28453.59
